$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to append after the last existing data row (row 95 -> new row 96)
$newRow = 96

$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value = "Maule"

# Date value (serial 44911 = 2022-12-16), keep same date style as the cell above (D95)
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item(95, 4).NumberFormat
$ws.Cells.Item($newRow, 4).Value = 44911

$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100101
$ws.Cells.Item($newRow, 8).Value = "Berries"
$ws.Cells.Item($newRow, 9).Value = 100101001
$ws.Cells.Item($newRow, 10).Value = "Arándano (blue)"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 220
$ws.Cells.Item($newRow, 14).Value = 3000
$ws.Cells.Item($newRow, 15).Value = 3000
$ws.Cells.Item($newRow, 16).Value = 3000
$ws.Cells.Item($newRow, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item($newRow, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($newRow, 19).Value = 1500
$ws.Cells.Item($newRow, 20).Value = 2
